# "files with surveyor names"
# Insert a new "surveyor" column (E) before the existing "notes" column,
# which shifts notes to column F. Fill in the surveyor name for every data
# row, then update the selection to match the recorded post-edit state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing "notes" values (currently in column E) before we
# overwrite that column with surveyor data, then re-write them into the new
# column F.
$notes = @{}
for ($r = 2; $r -le 23; $r++) {
    $v = $ws.Cells.Item($r, 5).Value2
    if ($v -ne $null -and $v -ne "") {
        $notes[$r] = $v
    }
}

# Shift notes from column E to column F.
for ($r = 2; $r -le 23; $r++) {
    if ($notes.ContainsKey($r)) {
        $ws.Cells.Item($r, 6).Value = $notes[$r]
    }
}

# New header for column E, and re-set the (shifted) header for column F.
$ws.Cells.Item(1, 5).Value = "surveyor"
$ws.Cells.Item(1, 6).Value = "notes"

# Surveyor name per row (bag_num rows 2-23).
$surveyors = @{
    2  = "Ava"
    3  = "Megan"
    4  = "Grace"
    5  = "Grace"
    6  = "Grace"
    7  = "Grace"
    8  = "Megan"
    9  = "Megan"
    10 = "Ava"
    11 = "Megan"
    12 = "Ava"
    13 = "Ava"
    14 = "Grace"
    15 = "Grace"
    16 = "Grace"
    17 = "Grace"
    18 = "Grace"
    19 = "Megan"
    20 = "Grace"
    21 = "Ava"
    22 = "Megan"
    23 = "Ava"
}

foreach ($r in $surveyors.Keys) {
    $ws.Cells.Item($r, 5).Value = $surveyors[$r]
}

# Update the active selection to match the saved workbook state.
$ws.Range("G27").Select()
